# Adds new plot/weather/aphid-count columns (mean_temp_year, total_rainfall_year,
# total_Metopolophium dirhodum, total_Rhopalosiphum padi, total_Sitobion avenae)
# to the "Mean" and "StdError" sheets, shifting the former "strip"/"CO2e_total_kg"
# layout (H:I) out to H:M, renaming H1 -> mean_temp_year, and re-homing the old
# CO2e_total_kg values into the new column M.

$wb = $excel.ActiveWorkbook

# Per-row data for the "Mean" sheet (columns H..M).
$meanData = @(
    @{ row=2;  H="10.4713470515";      I="709.0035599999999";  J="579";  K="1256"; L="324";   M="1699.532" },
    @{ row=3;  H="9.14304609925";      I="776.3224200000001";  J="415";  K="2683"; L="1115";  M="1699.532" },
    @{ row=4;  H="10.002850274";       I="714.42002";          J="4775"; K="4415"; L="4911";  M="1699.532" },
    @{ row=5;  H="9.48823883725";      I="935.7211";           J="381";  K="1349"; L="2123";  M="1699.532" },
    @{ row=6;  H="9.5933518805";       I="934.50516";          J="2583"; K="4376"; L="6747";  M="1699.532" },
    @{ row=7;  H="10.57316965825";     I="706.9033000000002";  J="229";  K="820";  L="3983";  M="1699.532" },
    @{ row=8;  H="9.243593227250001";  I="652.7387";           J="1101"; K="2169"; L="12803"; M="1699.532" },
    @{ row=9;  H="9.92623058975";      I="633.1731199999999";  J="5640"; K="4720"; L="2247";  M="1699.532" },
    @{ row=10; H="10.25889257125";     I="726.35834";          J="41";   K="1374"; L="385";   M="1699.532" },
    @{ row=11; H="10.29818976425";     I="951.63886";          J="265";  K="1757"; L="920";   M="1699.532" },
    @{ row=12; H="10.263630256";       I="794.0088199999999";  J="75";   K="1351"; L="877";   M="1699.532" },
    @{ row=13; H="9.70129622";         I="717.3";              J="2448"; K="2040"; L="2382";  M="1664.532" },
    @{ row=14; H="9.629243245750001";  I="623.7";              J="2519"; K="4968"; L="5991";  M="1664.532" },
    @{ row=15; H="10.525347886";       I="870.8";              J="86";   K="1833"; L="434";   M="1664.532" },
    @{ row=16; H="9.14851651625";      I="778.3000000000001";  J="151";  K="1081"; L="437";   M="1664.532" },
    @{ row=17; H="10.87966458975";     I="982.9";              J="150";  K="4315"; L="384";   M="1664.532" },
    @{ row=18; H="10.7708836005";      I="742.8";              J="110";  K="927";  L="111";   M="1664.532" },
    @{ row=19; H="10.6930498805";      I="622.2";              J="91";   K="660";  L="106";   M="1664.532" },
    @{ row=20; H="10.5555822335";      I="704.1999999999999";  J="191";  K="3848"; L="264";   M="1664.532" },
    @{ row=21; H="10.23620925";        I="565.2";              J="200";  K="843";  L="145";   M="1664.532" }
)

# Per-row data for the "StdError" sheet (columns I..M). H keeps its existing value (0).
$stdErrData = @(
    @{ row=2;  J="0"; K="0"; L="0"; M="0" },
    @{ row=3;  J="0"; K="0"; L="0"; M="0" },
    @{ row=4;  J="0"; K="0"; L="0"; M="0" },
    @{ row=5;  J="0"; K="0"; L="0"; M="0" },
    @{ row=6;  I="4.641245714809236e-14"; J="0"; K="0"; L="0"; M="9.282491429618472e-14" },
    @{ row=7;  J="0"; K="0"; L="0"; M="0" },
    @{ row=8;  J="0"; K="0"; L="0"; M="0" },
    @{ row=9;  J="0"; K="0"; L="0"; M="0" },
    @{ row=10; J="0"; K="0"; L="0"; M="0" },
    @{ row=11; J="0"; K="0"; L="0"; M="0" },
    @{ row=12; J="0"; K="0"; L="0"; M="0" },
    @{ row=13; J="0"; K="0"; L="0"; M="0" },
    @{ row=14; J="0"; K="0"; L="0"; M="0" },
    @{ row=15; J="0"; K="0"; L="0"; M="0" },
    @{ row=16; J="0"; K="0"; L="0"; M="0" },
    @{ row=17; J="0"; K="0"; L="0"; M="0" },
    @{ row=18; I="4.641245714809236e-14"; J="0"; K="0"; L="0"; M="9.282491429618472e-14" },
    @{ row=19; J="0"; K="0"; L="0"; M="0" },
    @{ row=20; J="0"; K="0"; L="0"; M="0" },
    @{ row=21; I="0"; J="0"; K="0"; L="0"; M="0" }
)

function Update-YieldSheet($sheetName, $rowData) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Stretch the header formatting (bold / bordered / centered, same as the
    # existing "strip" / "CO2e_total_kg" header cells) across the new columns.
    $ws.Range("I1").Copy()
    $ws.Range("J1:M1").PasteSpecial(-4122)

    # Rename / add header labels.
    $ws.Range("H1").Value = "mean_temp_year"
    $ws.Range("I1").Value = "total_rainfall_year"
    $ws.Range("J1").Value = "total_Metopolophium dirhodum"
    $ws.Range("K1").Value = "total_Rhopalosiphum padi"
    $ws.Range("L1").Value = "total_Sitobion avenae"
    $ws.Range("M1").Value = "CO2e_total_kg"

    foreach ($d in $rowData) {
        $r = $d.row
        if ($d.ContainsKey("H")) {
            $ws.Cells.Item($r, 8).Value = [double]$d.H
        }
        if ($d.ContainsKey("I")) {
            $ws.Cells.Item($r, 9).Value = [double]$d.I
        }
        $ws.Cells.Item($r, 10).Value = [double]$d.J
        $ws.Cells.Item($r, 11).Value = [double]$d.K
        $ws.Cells.Item($r, 12).Value = [double]$d.L
        $ws.Cells.Item($r, 13).Value = [double]$d.M
    }
}

Update-YieldSheet "Mean" $meanData
Update-YieldSheet "StdError" $stdErrData
